$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NroSiniestro (claim number) for the "preproduccion" row (row 3)
# from 1120194100378 to 1120170200917 (with two trailing spaces, kept as text)
$ws.Range("E3").Formula = "'1120170200917  "

# Update the active selection to L7 (was L10)
$ws.Range("L7").Select()
